# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7233
$ws1.Range("F4").Value = 122
$ws1.Range("F5").Value = 175
$ws1.Range("F7").Value = 94
$ws1.Range("F8").Value = 608

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7233
$ws4.Range("F5").Value = 122
$ws4.Range("F6").Value = 175
$ws4.Range("F9").Value = 94
$ws4.Range("F10").Value = 608
